$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-07 Friday" "2025-02-08 Saturday"

Replace-Text "284×2=" "128×3="
Replace-Text "577×2=" "452×7="
Replace-Text "179×2=" "516×7="
Replace-Text "448×2=" "646×8="
Replace-Text "216×3=" "527×3="

Replace-Text "910×9=" "337×8="
Replace-Text "192×8=" "813×4="
Replace-Text "484×6=" "139×5="
Replace-Text "443×4=" "845×4="
Replace-Text "424×2=" "440×9="

Replace-Text "862×3=" "488×9="
Replace-Text "274×3=" "342×3="
Replace-Text "785×4=" "121×5="
Replace-Text "268×5=" "867×2="
Replace-Text "102×3=" "548×2="

Replace-Text "895×7=" "867×4="
Replace-Text "913×2=" "409×8="
Replace-Text "129×4=" "186×7="
Replace-Text "941×5=" "747×3="
Replace-Text "705×5=" "694×5="

Replace-Text "390×6=" "658×7="
Replace-Text "375×2=" "873×2="
Replace-Text "520×2=" "745×6="
Replace-Text "479×6=" "942×8="
Replace-Text "331×3=" "303×2="

Write-Output "Done"
